# Adds two new FTT-Power variables (rows) used for the rooftop-solar /
# utility-electricity market-share calculation:
#   MEWDH - FTT-Power household demand
#   PRICH - FTT-Price of electricity use (incl taxes)
#
# Each variable gets one row on the "FTT-P" variable-definition sheet and
# one matching row on the "Time_Horizons" sheet.

$wb  = $excel.ActiveWorkbook
$wsP = $wb.Worksheets.Item("FTT-P")
$wsT = $wb.Worksheets.Item("Time_Horizons")

# --- FTT-P: new variable MEWDH (row 19) ------------------------------
$wsP.Range("A19").Value = "MEWDH"
$wsP.Range("B19").Value = 1
$wsP.Range("C19").Value = 33333333
$wsP.Range("D19").Value = "FTT-Power household demand"
$wsP.Range("E19").Value = "RSHORTTI"
$wsP.Range("F19").Value = "TIME"
$wsP.Range("G19").Value = 0
$wsP.Range("H19").Value = 0
$wsP.Range("I19").Value = "All"

# --- Time_Horizons: MEWDH time horizon (row 68) -----------------------
$wsT.Range("A68").Value = "MEWDH"
$wsT.Range("B68").Value = "tl_2010"

# --- FTT-P: new variable PRICH (row 20) -------------------------------
$wsP.Range("A20").Value = "PRICH"
$wsP.Range("B20").Value = 1
$wsP.Range("C20").Value = 33333334
$wsP.Range("D20").Value = "FTT-Price of electricity use (incl taxes)"
$wsP.Range("E20").Value = "RSHORTTI"
$wsP.Range("F20").Value = "TIME"
$wsP.Range("G20").Value = 0
$wsP.Range("H20").Value = 0
$wsP.Range("I20").Value = "All"

# --- Time_Horizons: PRICH time horizon (row 69) -----------------------
$wsT.Range("A69").Value = "PRICH"
$wsT.Range("B69").Value = "tl_2010"
